# Updates cryptos list (price/volume refresh + OKB/Mantle/ONDO reorder).
# All Price/Volume/Coin/Link columns are text (inline strings), even when
# their content looks numeric (e.g. "1.00", "0.999", "7.60"). Plain
# Value2 assignment is enough for values Excel can't parse as a plain
# number (contain two dots, %, spaces, letters, ...). For values that
# ARE valid numeric literals, force the cell to Text ("@") first so the
# literal text (with its exact trailing zeros) is preserved instead of
# being coerced into a number, then drop back to the default "Normal"
# style so no stray number-format style is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "61.314.33"
$ws.Range("E2").Value2 = "  +0.07%  "
$ws.Range("D3").Value2 = "3.374.83"
$ws.Range("E3").Value2 = "  +1.90%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value2 = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "573.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "  +1.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "136.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = "  +6.70%  "
$ws.Range("E7").Value2 = "  -0.11%  "
$ws.Range("D8").Value2 = "3.377.23"
$ws.Range("E8").Value2 = "  +2.02%  "
$ws.Range("E9").Value2 = "  -0.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "7.60"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value2 = "  +3.47%  "
$ws.Range("E11").Value2 = "  +4.33%  "
$ws.Range("E12").Value2 = "  +3.81%  "
$ws.Range("D13").Value2 = "3.942.22"
$ws.Range("E13").Value2 = "  +1.60%  "
$ws.Range("E14").Value2 = "  +2.00%  "
$ws.Range("E15").Value2 = "  +4.74%  "
$ws.Range("D16").Value2 = "3.372.47"
$ws.Range("E16").Value2 = "  +1.81%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "25.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value2 = "  +1.95%  "
$ws.Range("D18").Value2 = "61.296.52"
$ws.Range("E18").Value2 = "  -0.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "5.90"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = "  +3.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "13.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = "  +2.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "9.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value2 = "  +3.80%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "380.95"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value2 = "  +7.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "0.568"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value2 = "  +2.18%  "
$ws.Range("D24").Value2 = "3.499.85"
$ws.Range("E24").Value2 = "  +1.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value2 = "  +0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "70.70"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value2 = "  +1.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "0.0000120"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value2 = "  +11.26%  "
$ws.Range("E28").Value2 = "  +13.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "7.77"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value2 = "  +8.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "0.998"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value2 = "  -0.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "8.22"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value2 = "  +4.52%  "
$ws.Range("E32").Value2 = "  +3.93%  "
$ws.Range("E33").Value2 = "  +1.21%  "
$ws.Range("E34").Value2 = "  -0.10%  "
$ws.Range("D35").Value2 = "3.395.75"
$ws.Range("E35").Value2 = "  +1.56%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "23.46"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value2 = "  +3.90%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "5.50"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value2 = "  +4.72%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "7.04"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value2 = "  +3.14%  "
$ws.Range("E39").Value2 = "  +4.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "160.52"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value2 = "  +0.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "0.0788"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value2 = "  +3.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value2 = "  -0.13%  "
$ws.Range("E43").Value2 = "  +10.56%  "
$ws.Range("E44").Value2 = "  +1.00%  "
$ws.Range("B45").Value2 = "ONDO"
$ws.Range("C45").Value2 = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "1.21"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value2 = "  +8.70%  "
$ws.Range("B46").Value2 = "OKB"
$ws.Range("C46").Value2 = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "41.55"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value2 = "  +1.29%  "
$ws.Range("B47").Value2 = "Mantle"
$ws.Range("C47").Value2 = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "0.764"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value2 = "  +2.73%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "23.06"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value2 = "  +3.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "6.97"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value2 = "  +3.67%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "22.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value2 = "  +7.52%  "
$ws.Range("D51").Value2 = "2.322.01"
$ws.Range("E51").Value2 = "  +7.64%  "
